# Updated cryptos list on Sat Apr 15 20:28:53 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures for each
# coin row (rows 2-51) on the active worksheet to match the latest scrape.
#
# Column D prices are stored as plain text (e.g. "30.455.11", "1.006") so
# that multi-dot thousands separators survive untouched. Values that look
# numeric (e.g. "1.006", "0.08920") would otherwise be silently converted to
# real numbers by Excel on assignment through .Value - dropping meaningful
# trailing zeros / formatting. To avoid that, such cells are briefly switched
# to Text number format before the write and then have that formatting
# cleared again right after, leaving the cell's visual style untouched while
# keeping the stored value as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cell, $text, $forceText) {
    if ($forceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

# Row 2
Set-PriceText $ws.Range("D2") "30.455.11" $false
$ws.Range("E2").Value = "  +0.34%  "

# Row 3
Set-PriceText $ws.Range("D3") "2.107.46" $false
$ws.Range("E3").Value = "  +0.62%  "

# Row 4
Set-PriceText $ws.Range("D4") "1.006" $true
$ws.Range("E4").Value = "  +0.63%  "

# Row 5
Set-PriceText $ws.Range("D5") "334.57" $true
$ws.Range("E5").Value = "  +1.73%  "

# Row 6
Set-PriceText $ws.Range("D6") "1.004" $true
$ws.Range("E6").Value = "  +0.46%  "

# Row 7
Set-PriceText $ws.Range("D7") "0.5217" $true
$ws.Range("E7").Value = "  -0.24%  "

# Row 8
Set-PriceText $ws.Range("D8") "0.4531" $true
$ws.Range("E8").Value = "  +4.32%  "

# Row 9
Set-PriceText $ws.Range("D9") "54.47" $true
$ws.Range("E9").Value = "  +16.33%  "

# Row 10
Set-PriceText $ws.Range("D10") "0.08920" $true
$ws.Range("E10").Value = "  +1.10%  "

# Row 11
Set-PriceText $ws.Range("D11") "1.181" $true
$ws.Range("E11").Value = "  +1.44%  "

# Row 12
Set-PriceText $ws.Range("D12") "24.13" $true
$ws.Range("E12").Value = "  -1.48%  "

# Row 13
Set-PriceText $ws.Range("D13") "2.102.66" $false
$ws.Range("E13").Value = "  +0.69%  "

# Row 14
Set-PriceText $ws.Range("D14") "6.820" $true
$ws.Range("E14").Value = "  +1.36%  "

# Row 15
Set-PriceText $ws.Range("D15") "8.014" $true
$ws.Range("E15").Value = "  +3.37%  "

# Row 16
Set-PriceText $ws.Range("D16") "96.72" $true
$ws.Range("E16").Value = "  +0.21%  "

# Row 17
$ws.Range("E17").Value = "  +1.30%  "

# Row 18
$ws.Range("E18").Value = "  +0.44%  "

# Row 19
Set-PriceText $ws.Range("D19") "0.06648" $true
$ws.Range("E19").Value = "  +0.15%  "

# Row 20
Set-PriceText $ws.Range("D20") "19.21" $true
$ws.Range("E20").Value = "  +1.36%  "

# Row 21
$ws.Range("E21").Value = "  +0.46%  "

# Row 22
Set-PriceText $ws.Range("D22") "6.327" $true
$ws.Range("E22").Value = "  -0.29%  "

# Row 23
Set-PriceText $ws.Range("D23") "30.509.75" $false
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
Set-PriceText $ws.Range("D24") "12.41" $true
$ws.Range("E24").Value = "  +0.33%  "

# Row 25
Set-PriceText $ws.Range("D25") "2.350" $true
$ws.Range("E25").Value = "  +1.71%  "

# Row 26
Set-PriceText $ws.Range("D26") "2.350.89" $false
$ws.Range("E26").Value = "  +0.74%  "

# Row 27
Set-PriceText $ws.Range("D27") "22.16" $true
$ws.Range("E27").Value = "  -1.11%  "

# Row 28
Set-PriceText $ws.Range("D28") "162.72" $true
$ws.Range("E28").Value = "  +0.60%  "

# Row 29
Set-PriceText $ws.Range("D29") "2.530" $true
$ws.Range("E29").Value = "  -2.51%  "

# Row 30
Set-PriceText $ws.Range("D30") "133.94" $true
$ws.Range("E30").Value = "  +1.32%  "

# Row 31
Set-PriceText $ws.Range("D31") "1.208" $true
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
Set-PriceText $ws.Range("D32") "0.1069" $true
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
Set-PriceText $ws.Range("D33") "6.386" $true
$ws.Range("E33").Value = "  +3.44%  "

# Row 34
Set-PriceText $ws.Range("D34") "1.637" $true
$ws.Range("E34").Value = "  -2.81%  "

# Row 35
Set-PriceText $ws.Range("D35") "3.943" $true
$ws.Range("E35").Value = "  +1.14%  "

# Row 36
Set-PriceText $ws.Range("D36") "10.36" $true
$ws.Range("E36").Value = "  +3.60%  "

# Row 37
Set-PriceText $ws.Range("D37") "5.768" $true
$ws.Range("E37").Value = "  +5.29%  "

# Row 38
Set-PriceText $ws.Range("D38") "0.02583" $true
$ws.Range("E38").Value = "  -0.07%  "

# Row 39
Set-PriceText $ws.Range("D39") "0.06828" $true
$ws.Range("E39").Value = "  +1.91%  "

# Row 40
Set-PriceText $ws.Range("D40") "0.2303" $true
$ws.Range("E40").Value = "  +1.95%  "

# Row 41
Set-PriceText $ws.Range("D41") "12.73" $true
$ws.Range("E41").Value = "  +0.48%  "

# Row 42
Set-PriceText $ws.Range("D42") "0.6867" $true
$ws.Range("E42").Value = "  +0.65%  "

# Row 43
Set-PriceText $ws.Range("D43") "1.246" $true
$ws.Range("E43").Value = "  -0.22%  "

# Row 44
$ws.Range("E44").Value = "  +5.03%  "

# Row 45
Set-PriceText $ws.Range("D45") "14.08" $true
$ws.Range("E45").Value = "  -0.12%  "

# Row 46
Set-PriceText $ws.Range("D46") "0.6359" $true
$ws.Range("E46").Value = "  -0.21%  "

# Row 47
Set-PriceText $ws.Range("D47") "3.663" $true
$ws.Range("E47").Value = "  +1.26%  "

# Row 48
$ws.Range("E48").Value = "  +21.56%  "

# Row 49
$ws.Range("E49").Value = "  +0.03%  "

# Row 50
Set-PriceText $ws.Range("D50") "1.203" $true
$ws.Range("E50").Value = "  +0.86%  "

# Row 51
Set-PriceText $ws.Range("D51") "83.07" $true
$ws.Range("E51").Value = "  +1.39%  "

